{"js": "const replacements = [\n  [\"2024-04-25 Thursday\", \"2024-04-26 Friday\"],\n  [\"24\u00f78=\", \"70\u00f74=\"],\n  [\"90\u00f78=\", \"53\u00f74=\"],\n  [\"82\u00f78=\", \"14\u00f78=\"],\n  [\"58\u00f79=\", \"52\u00f76=\"],\n  [\"22\u00f73=\", \"21\u00f77=\"],\n  [\"68\u00f74=\", \"28\u00f72=\"],\n  [\"95\u00f78=\", \"91\u00f77=\"],\n  [\"70\u00f78=\", \"98\u00f77=\"],\n  [\"21\u00f72=\", \"46\u00f78=\"],\n  [\"56\u00f78=\", \"30\u00f76=\"],\n  [\"81\u00f76=\", \"23\u00f78=\"],\n  [\"92\u00f73=\", \"52\u00f75=\"],\n  [\"21\u00f78=\", \"57\u00f78=\"],\n  [\"32\u00f75=\", \"87\u00f78=\"],\n  [\"43\u00f76=\", \"18\u00f73=\"],\n  [\"10\u00f78=\", \"45\u00f78=\"],\n  [\"70\u00f79=\", \"78\u00f72=\"],\n  [\"75\u00f79=\", \"78\u00f74=\"],\n  [\"12\u00f74=\", \"12\u00f73=\"],\n  [\"51\u00f79=\", \"58\u00f74=\"],\n  [\"34\u00f72=\", \"57\u00f78=\"],\n  [\"11\u00f74=\", \"92\u00f73=\"],\n  [\"97\u00f77=\", \"85\u00f73=\"],\n  [\"94\u00f72=\", \"55\u00f79=\"],\n  [\"80\u00f79=\", \"33\u00f75=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-25 Thursday\", \"2024-04-26 Friday\"),\n    @(\"24\u00f78=\", \"70\u00f74=\"),\n    @(\"90\u00f78=\", \"53\u00f74=\"),\n    @(\"82\u00f78=\", \"14\u00f78=\"),\n    @(\"58\u00f79=\", \"52\u00f76=\"),\n    @(\"22\u00f73=\", \"21\u00f77=\"),\n    @(\"68\u00f74=\", \"28\u00f72=\"),\n    @(\"95\u00f78=\", \"91\u00f77=\"),\n    @(\"70\u00f78=\", \"98\u00f77=\"),\n    @(\"21\u00f72=\", \"46\u00f78=\"),\n    @(\"56\u00f78=\", \"30\u00f76=\"),\n    @(\"81\u00f76=\", \"23\u00f78=\"),\n    @(\"92\u00f73=\", \"52\u00f75=\"),\n    @(\"21\u00f78=\", \"57\u00f78=\"),\n    @(\"32\u00f75=\", \"87\u00f78=\"),\n    @(\"43\u00f76=\", \"18\u00f73=\"),\n    @(\"10\u00f78=\", \"45\u00f78=\"),\n    @(\"70\u00f79=\", \"78\u00f72=\"),\n    @(\"75\u00f79=\", \"78\u00f74=\"),\n    @(\"12\u00f74=\", \"12\u00f73=\"),\n    @(\"51\u00f79=\", \"58\u00f74=\"),\n    @(\"34\u00f72=\", \"57\u00f78=\"),\n    @(\"11\u00f74=\", \"92\u00f73=\"),\n    @(\"97\u00f77=\", \"85\u00f73=\"),\n    @(\"94\u00f72=\", \"55\u00f79=\"),\n    @(\"80\u00f79=\", \"33\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($from, $false, $false, $false, $false, $false, $true, 1, $false, $to, 2)\n}\n"}
